$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Long text values (shared strings) ---
$objetivosVal = "Fornecer aos alunos capacidade para: Decidir sobre a melhor alternativa a ser adotada para tratamento avançado, em função da qualidade do efluente a tratar e dos objetivos do reuso ou limitações do corpo receptor. Conhecer as tecnologias disponíveis para tratamento avançado de águas residuárias. Dispor de conceitos e conhecimentos para pré-dimensionamento e especificação de equipamentos para sistemas de tratamento avançado."
$docenteVal = "7455355 - Robson da Silva Rocha"
$programaResumidoVal = "Tendências mundiais no tratamento de águas. Processo e operação na remoção de elementos indesejáveis. Recuperação de Ambientes Aquáticos."
$programaVal = "Tendências mundiais sobre tratamento avançado e reuso de águas residuárias. Determinação da eficiência de processos e operações em função dos objetivos de reuso da qualidade do afluente a tratar e da obediência a padrões de emissão e de qualidade. Processos e operações aplicadas à remoção de nitrogênio e fósforo: nitrificação, desnitrificação, e remoção química e biológica de fósforo. Operações e processos para remoção de contaminantes específicos: adsorção em carvão ativado, oxidação química, `"stripping`", coagulação-floculação (sedimentação e flotação), troca iônica, osmose reversa, filtração em membranas, filtração em meios porosos. Disposição e tratamento de esgotos no solo. Recuperação de ambientes aquáticos com base na piscicultura e aproveitamento de algas e macrófias."
$biblioVal = "CHERNICHARO, C.A.L. (coord.). Pós-Tratamento de Efluentes de Reatores Anaeróbios - PROCOPE - PROSAB - ISBN - 85 - 901640-1-2, Belo Horizonte, MG, 220p.`nCAMPOS, J.R. (1994), Avaliação do Processo Eletrolítico para Tratamento de Esgoto Sanitário. Consórcio Intermunicipal das Bacias dos Rios Piracicaba e Capivari, Americana, SP, 12p.`nCAMPOS, J.R. (1994). Pré-tratamento de Águas para Abastecimento. ASSEMAE - Consórcio Intermunicipal das Bacias dos Rios Piracicaba e Capivari. Publ. n. 9, Americana, SP, 112p.`nDANIEL, L.A.; CAMPOS, J.R. (1993). Radiação Ultravioleta é Alternativa Viável para Desinfecção de Efluentes de Sistemas de Tratamento Aeróbio e Anaeróbio no Brasil, BIO, n. 5, set/out, p. 1-17.`nLAPOLLI, F.R. (1998). Processos de Separação por Membranas. São Carlos - SP, 76p. (apostila)."

# --- Row 10: Objetivos answer text replaced ---
$ws.Range("B10").Value2 = $objetivosVal
$ws.Range("C10").Value2 = $objetivosVal

# --- Row 13: now holds the "Docentes responsaveis" value, label cell cleared ---
$ws.Range("A13").Value2 = ""
$ws.Range("B13").Value2 = $docenteVal
$ws.Range("C13").Value2 = $docenteVal

# --- Row 14: "Programa resumido:" label + new summary text ---
$ws.Range("A14").Value2 = "Programa resumido:"
$ws.Range("B14").Value2 = $programaResumidoVal
$ws.Range("C14").Value2 = $programaResumidoVal

# --- Row 15: "Short syllabus:" label, value cells cleared ---
$ws.Range("A15").Value2 = "Short syllabus:"
$ws.Range("B15").Value2 = ""
$ws.Range("C15").Value2 = ""

# --- Row 16: "Programa:" label + new programa text ---
$ws.Range("A16").Value2 = "Programa:"
$ws.Range("B16").Value2 = $programaVal
$ws.Range("C16").Value2 = $programaVal

# --- Row 17: "Syllabus:" label only ---
$ws.Range("A17").Value2 = "Syllabus:"

# --- Row 18: "Avaliacao:" label, value cells cleared ---
$ws.Range("A18").Value2 = "Avaliação:"
$ws.Range("B18").Value2 = ""
$ws.Range("C18").Value2 = ""

# --- Row 19: "Metodo:" label (text value already correct) ---
$ws.Range("A19").Value2 = "Método:"

# --- Row 20: "Criterio:" label (text value already correct) ---
$ws.Range("A20").Value2 = "Critério:"

# --- Row 21: "Norma de recuperacao:" label (text value already correct) ---
$ws.Range("A21").Value2 = "Norma de recuperação:"

# --- Row 22 (new row): "Bibliografia:" label + bibliography text ---
$ws.Range("A22").Value2 = "Bibliografia:"
$ws.Range("B22").Value2 = $biblioVal
$ws.Range("C22").Value2 = $biblioVal

# --- Row heights ---
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 15
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120

# --- Column width fix: column A keeps 30.71, column B becomes its own 60.71 definition ---
$ws.Range("A1").EntireColumn.ColumnWidth = 30.7109375
$ws.Range("B1").EntireColumn.ColumnWidth = 60.7109375

